# Sketched a round UI. Added all of the inventory elements to a grid layout
#
# Inserts a new "Controls:" section (with a blank separator paragraph)
# right after the existing "Dodge Chance" paragraph and before the
# trailing (whitespace-only) paragraph that precedes the section break.

$d = $word.ActiveDocument
$wordml = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# The paragraphs to insert, in document order, expressed as raw WordprocessingML.
$newParagraphs = @(
    "<w:p $wordml/>",
    "<w:p $wordml><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Controls:</w:t></w:r></w:p>",
    "<w:p $wordml><w:r><w:t>A: attack</w:t></w:r></w:p>",
    "<w:p $wordml><w:r><w:t>X and Y: attack buttons</w:t></w:r></w:p>",
    "<w:p $wordml><w:r><w:t>4 hotbar slots: dpad</w:t></w:r></w:p>",
    "<w:p $wordml><w:r><w:lastRenderedPageBreak/><w:t xml:space=`"preserve`">Up to 4 abilities: </w:t></w:r><w:r><w:t>c stick</w:t></w:r></w:p>",
    "<w:p $wordml><w:r><w:t>Use alternate weapon: zl</w:t></w:r></w:p>",
    "<w:p $wordml><w:r><w:t>Dash: L and R</w:t></w:r></w:p>"
)

# Anchor on the start of the final (whitespace-only) paragraph that sits
# right before the section break - every new paragraph is inserted there,
# in order, which pushes that anchor paragraph further down each time.
$lastParaIndex = $d.Paragraphs.Count
$anchorStart = $d.Paragraphs.Item($lastParaIndex).Range.Start

foreach ($xml in $newParagraphs) {
    $insertionPoint = $d.Range($anchorStart, $anchorStart)
    $insertionPoint.InsertParagraphBefore()

    $lastParaIndex = $d.Paragraphs.Count
    $newParaRange = $d.Paragraphs.Item($lastParaIndex - 1).Range
    $newParaRange.InsertXML($xml) | Out-Null

    $anchorStart = $d.Paragraphs.Item($d.Paragraphs.Count).Range.Start
}

Write-Output "Inserted $($newParagraphs.Count) paragraphs; document now has $($d.Paragraphs.Count) paragraphs."
